# Processors: Performance: Overclocking: Add VBS off
#
# Inserts 3 new rows (16-18) for a new "PBO Enhanced 3 + VBS off" profile,
# pushing the existing "PBO Enhanced 3 + Perf Switch 1/2" and "9900X" rows
# down by three (old rows 16-22 -> new rows 19-25).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert three blank rows before the old row 16 -----------------
$ws.Rows("16:18").Insert()

# Row-insert copies the per-cell format of the row above (row 15, which
# has data in O/Q) into the new rows even though this profile's rows
# never use those columns. Drop the stray formatted-but-empty cells so
# rows 16-18 only contain D:J, matching the rest of this profile block.
$ws.Range("O16:O18").Clear()
$ws.Range("Q16:Q18").Clear()

# --- 2. Populate the new rows with the "PBO Enhanced 3 + VBS off" data -
$ws.Range("D16").Value2 = "PBO Enhanced 3 + VBS off"
$ws.Range("E16").Value2 = 4443
$ws.Range("F16").Formula = "=E16/`$E`$3"
$ws.Range("G16").Value2 = 55650
$ws.Range("H16").Formula = "=G16/`$G`$3"
$ws.Range("I16").Formula = "=G16/E16"
$ws.Range("J16").Formula = "=G16/12"

$ws.Range("D17").Value2 = "PBO Enhanced 3 + VBS off"
$ws.Range("E17").Value2 = 4438
$ws.Range("F17").Formula = "=E17/`$E`$3"
$ws.Range("G17").Value2 = 55739
$ws.Range("H17").Formula = "=G17/`$G`$3"
$ws.Range("I17").Formula = "=G17/E17"
$ws.Range("J17").Formula = "=G17/12"

$ws.Range("D18").Value2 = "PBO Enhanced 3 + VBS off"
$ws.Range("E18").Value2 = 4439
$ws.Range("F18").Formula = "=E18/`$E`$3"
$ws.Range("G18").Value2 = 55610
$ws.Range("H18").Formula = "=G18/`$G`$3"
$ws.Range("I18").Formula = "=G18/E18"
$ws.Range("J18").Formula = "=G18/12"

# --- 3. Selection / cursor position, as recorded in the saved file -----
$ws.Range("K20").Select()
